# Daily update of covid19 tracker data files
# Bump the "Last Updated" date in column B (rows 5-74) of the
# "Country Updates" sheet from 4/4/2020 (43925) to 4/5/2020 (43926).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Country Updates")

$oldDate = [DateTime]::FromOADate(43925)
$newDate = [DateTime]::FromOADate(43926)

for ($r = 5; $r -le 74; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $val = $cell.Value()
    if ($val -ne $null -and $val.ToOADate() -eq 43925) {
        $cell.Value = $newDate
    }
}
